$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Constant text used across sheets
# -----------------------------------------------------------------------
$mdNew      = "72e708ce-bc17-4331-b2fb-5662bc36f88f.md"
$cfgName    = ".localization-config"
$png1       = "72fb6fd1-ca71-4751-b653-6552835a64ea.png"
$png2       = "c0263322-3c25-480c-b8aa-11908b6fc81a.png"
$zhcnXlf    = "72e708ce-bc17-4331-b2fb-5662bc36f88f.87a0184ba93c0178eb3a3cb3786efca92aeb49d1.zh-cn.xlf"
$dedeXlf    = "72e708ce-bc17-4331-b2fb-5662bc36f88f.87a0184ba93c0178eb3a3cb3786efca92aeb49d1.de-de.xlf"
$depPngA405 = "a405cace286a851630514a9c51bf0bfa3a33b50f.png"
$depPng3e37 = "3e37a7db278a398d910680fc70b59204594351bc.png"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$epoch           = "0001-01-01 00:00:00"
$dtZhCn          = "2016-03-10 07:33:40"
$dtDeDe          = "2016-03-10 07:33:46"
$include         = "Include"
$isDependency    = "IsDependency"
$ignored         = "Ignored"
$depMdPath       = "e2e\72e708ce-bc17-4331-b2fb-5662bc36f88f.md"

# Hyperlink URL bases (mirrors the existing convention already present
# in the workbook's relationships).
$mdBase     = "https://github.com/OpenLocalizationTest/oltest/blob/187c6c7e0fc8fb52ee696f4e5a477ee88f2decd2/e2e/"
$cfgBase    = "https://github.com/OpenLocalizationTest/oltest/blob/187c6c7e0fc8fb52ee696f4e5a477ee88f2decd2/"
$zhcnBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f0846764256162bca9b56b70c2601ec367ef68c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$dedeBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c0885fb7fcf113aaa0aebef0b641c10352ddb35/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

function Clear-AllHyperlinks($ws) {
    $existing = @()
    foreach ($h in $ws.Hyperlinks) { $existing += $h }
    for ($i = $existing.Length - 1; $i -ge 0; $i--) {
        $existing[$i].Delete()
    }
}

# -----------------------------------------------------------------------
# Sheet 1: Overview
# -----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

Clear-AllHyperlinks $wsOverview

$wsOverview.Range("A2").Value() = $mdNew
$wsOverview.Range("B2").Value() = $readyForHandoff
$wsOverview.Range("C2").Value() = $readyForHandoff

$wsOverview.Range("A3").Value() = $png1
$wsOverview.Range("B3").Value() = $readyForHandoff
$wsOverview.Range("C3").Value() = $readyForHandoff

$wsOverview.Range("A4").Value() = $png2
$wsOverview.Range("B4").Value() = $readyForHandoff
$wsOverview.Range("C4").Value() = $readyForHandoff

$wsOverview.Range("A5").Value() = $cfgName
$wsOverview.Range("B5").Value() = $notLocalized
$wsOverview.Range("C5").Value() = $notLocalized

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$mdBase$mdNew", "", "", $mdNew) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$mdBase$png1", "", "", $png1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$mdBase$png2", "", "", $png2) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "$cfgBase$cfgName", "", "", $cfgName) | Out-Null

# -----------------------------------------------------------------------
# Sheet 2: zh-cn
# -----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

Clear-AllHyperlinks $wsZhCn

$wsZhCn.Range("A2").Value() = $mdNew
$wsZhCn.Range("B2").Value() = $readyForHandoff
$wsZhCn.Range("C2").Value() = $zhcnXlf
$wsZhCn.Range("D2").Value() = $dtZhCn
$wsZhCn.Range("G2").Value() = $epoch
$wsZhCn.Range("H2").Value() = $include

$wsZhCn.Range("A3").Value() = $png1
$wsZhCn.Range("B3").Value() = $readyForHandoff
$wsZhCn.Range("C3").Value() = $depPngA405
$wsZhCn.Range("D3").Value() = $dtZhCn
$wsZhCn.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("G3").Value() = $epoch
$wsZhCn.Range("H3").Value() = $isDependency
$wsZhCn.Range("I3").Value() = $depMdPath

$wsZhCn.Range("A4").Value() = $png2
$wsZhCn.Range("B4").Value() = $readyForHandoff
$wsZhCn.Range("C4").Value() = $depPng3e37
$wsZhCn.Range("D4").Value() = $dtZhCn
$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("G4").Value() = $epoch
$wsZhCn.Range("H4").Value() = $isDependency
$wsZhCn.Range("I4").Value() = $depMdPath

$wsZhCn.Range("A5").Value() = $cfgName
$wsZhCn.Range("B5").Value() = $notLocalized
$wsZhCn.Range("D5").Value() = $epoch
$wsZhCn.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("G5").Value() = $epoch
$wsZhCn.Range("H5").Value() = $ignored

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "$mdBase$mdNew", "", "", $mdNew) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "$zhcnBase$zhcnXlf", "", "", $zhcnXlf) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "$mdBase$png1", "", "", $png1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "$zhcnBase$depPngA405", "", "", $depPngA405) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "$mdBase$png2", "", "", $png2) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), "$zhcnBase$depPng3e37", "", "", $depPng3e37) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "$cfgBase$cfgName", "", "", $cfgName) | Out-Null

# -----------------------------------------------------------------------
# Sheet 3: de-de
# -----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

Clear-AllHyperlinks $wsDeDe

$wsDeDe.Range("A2").Value() = $mdNew
$wsDeDe.Range("B2").Value() = $readyForHandoff
$wsDeDe.Range("C2").Value() = $dedeXlf
$wsDeDe.Range("D2").Value() = $dtDeDe
$wsDeDe.Range("G2").Value() = $epoch
$wsDeDe.Range("H2").Value() = $include

$wsDeDe.Range("A3").Value() = $png1
$wsDeDe.Range("B3").Value() = $readyForHandoff
$wsDeDe.Range("C3").Value() = $depPngA405
$wsDeDe.Range("D3").Value() = $dtDeDe
$wsDeDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("G3").Value() = $epoch
$wsDeDe.Range("H3").Value() = $isDependency
$wsDeDe.Range("I3").Value() = $depMdPath

$wsDeDe.Range("A4").Value() = $png2
$wsDeDe.Range("B4").Value() = $readyForHandoff
$wsDeDe.Range("C4").Value() = $depPng3e37
$wsDeDe.Range("D4").Value() = $dtDeDe
$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("G4").Value() = $epoch
$wsDeDe.Range("H4").Value() = $isDependency
$wsDeDe.Range("I4").Value() = $depMdPath

$wsDeDe.Range("A5").Value() = $cfgName
$wsDeDe.Range("B5").Value() = $notLocalized
$wsDeDe.Range("D5").Value() = $epoch
$wsDeDe.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("G5").Value() = $epoch
$wsDeDe.Range("H5").Value() = $ignored

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "$mdBase$mdNew", "", "", $mdNew) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "$dedeBase$dedeXlf", "", "", $dedeXlf) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "$mdBase$png1", "", "", $png1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "$dedeBase$depPngA405", "", "", $depPngA405) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "$mdBase$png2", "", "", $png2) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), "$dedeBase$depPng3e37", "", "", $depPng3e37) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "$cfgBase$cfgName", "", "", $cfgName) | Out-Null

Write-Host "Report for handoff generated."
